$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("omp")

# New shared strings are interned in first-use order, so touch the text
# cells in the same order the target workbook assigned indices 26-28:
# "godel omp classifer non-blocked", "godel omp classifer blocked", "speedup".
$ws.Range("H8").Value = "godel omp classifer non-blocked"
$ws.Range("H14").Value = "godel omp classifer blocked"
$ws.Range("F2").Value = "speedup"

# --- New "speedup" column (F), computed from existing D/B columns ---
$ws.Range("F5").Formula = "=D5/B5"
$ws.Range("F6:F10").Formula = "=D6/B6"

# --- Non-blocked classifier micro-benchmark block (rows 8-11, cols H-P) ---
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 2
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 8
$ws.Range("L9").Value = 16
$ws.Range("M9").Value = 32
$ws.Range("N9").Value = 64
$ws.Range("O9").Value = 128
$ws.Range("P9").Value = 256

$ws.Range("H10").Value = 115
$ws.Range("I10").Value = 71
$ws.Range("J10").Value = 84
$ws.Range("K10").Value = 51
$ws.Range("L10").Value = 52
$ws.Range("M10").Value = 54
$ws.Range("N10").Value = 47
$ws.Range("O10").Value = 34
$ws.Range("P10").Value = 51

$ws.Range("H11").Formula = "=`$H`$10/H10"
$ws.Range("I11:P11").Formula = "=`$H`$10/I10"

# --- Blocked classifier micro-benchmark block (rows 14-17, cols H-P) ---
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 8
$ws.Range("L15").Value = 16
$ws.Range("M15").Value = 32
$ws.Range("N15").Value = 64
$ws.Range("O15").Value = 128
$ws.Range("P15").Value = 256

$ws.Range("H16").Value = 115
$ws.Range("I16").Value = 69
$ws.Range("J16").Value = 66
$ws.Range("K16").Value = 50
$ws.Range("L16").Value = 52
$ws.Range("M16").Value = 54
$ws.Range("N16").Value = 38
$ws.Range("O16").Value = 35
$ws.Range("P16").Value = 54

$ws.Range("H17").Formula = "=`$H`$16/H16"
$ws.Range("I17:P17").Formula = "=`$H`$16/I16"

# --- Selection moved to F3 (view-state only) ---
$ws.Activate() | Out-Null
$ws.Range("F3").Select() | Out-Null
